# The commit swaps the contents of ppt/theme/theme1.xml ("Office Theme"
# colours) and ppt/theme/theme2.xml ("Integral" colours). theme2.xml is the
# theme that backs the presentation's single Slide Master / Design, so the
# reachable, observable effect of that swap is that the design's twelve
# theme colours change from the "Integral" palette to the stock "Office
# Theme" palette (the font scheme and format scheme are already identical
# between the two themes, so nothing else visibly changes).
$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colorScheme = $theme.ThemeColorScheme

# Index -> (scheme slot, target "Office Theme" RGB as a COM long: R | G<<8 | B<<16)
$colorScheme.Item(1).RGB  = 0         # dk1      -> 000000
$colorScheme.Item(2).RGB  = 16777215  # lt1      -> FFFFFF
$colorScheme.Item(3).RGB  = 6968388   # dk2      -> 44546A
$colorScheme.Item(4).RGB  = 15132391  # lt2      -> E7E6E6
$colorScheme.Item(5).RGB  = 13998939  # accent1  -> 5B9BD5
$colorScheme.Item(6).RGB  = 3243501   # accent2  -> ED7D31
$colorScheme.Item(7).RGB  = 10855845  # accent3  -> A5A5A5
$colorScheme.Item(8).RGB  = 49407     # accent4  -> FFC000
$colorScheme.Item(9).RGB  = 12874308  # accent5  -> 4472C4
$colorScheme.Item(10).RGB = 4697456   # accent6  -> 70AD47
$colorScheme.Item(11).RGB = 12673797  # hlink    -> 0563C1
$colorScheme.Item(12).RGB = 7491477   # folHlink -> 954F72
